$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (machine_index = 0)
$ws.Range("B2").Value = 41
$ws.Range("C2").Value = 185
$ws.Range("D2").Value = 108
$ws.Range("E2").Value = 77
$ws.Range("F2").Value = 126
$ws.Range("G2").Value = 85.7

# Row 3 (machine_index = 1)
$ws.Range("B3").Value = 24
$ws.Range("C3").Value = 483
$ws.Range("E3").Value = 483

# Row 4 (machine_index = 2)
$ws.Range("B4").Value = 85
$ws.Range("C4").Value = 340
$ws.Range("D4").Value = 340
$ws.Range("F4").Value = 340

# Row 5 (machine_index = 3)
$ws.Range("B5").Value = 95
$ws.Range("C5").Value = 296
$ws.Range("D5").Value = 190
$ws.Range("E5").Value = 106
$ws.Range("F5").Value = 218
$ws.Range("G5").Value = 87.2

# Row 6 (machine_index = 4)
$ws.Range("C6").Value = 243
$ws.Range("E6").Value = 243
